$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.942.61'
$ws.Range("D3").Value = '2.666.89'
$ws.Range("E3").Value = '  -1.10%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.56'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '164.37'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.08%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.545'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.26%  '
$ws.Range("D9").Value = '2.664.44'
$ws.Range("E9").Value = '  -1.10%  '
$ws.Range("E10").Value = '  +1.17%  '
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.357'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.53%  '
$ws.Range("E13").Value = '  -1.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.68'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.42%  '
$ws.Range("D15").Value = '3.151.67'
$ws.Range("E15").Value = '  -1.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000182'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.54%  '
$ws.Range("D17").Value = '67.031.86'
$ws.Range("E17").Value = '  -2.04%  '
$ws.Range("D18").Value = '2.654.20'
$ws.Range("E18").Value = '  -1.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.61'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '360.25'
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.48'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.37'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.79'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.51%  '
$ws.Range("E24").Value = '  -5.24%  '
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '71.08'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -5.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.05'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.97%  '
$ws.Range("E28").Value = '  -1.42%  '
$ws.Range("E29").Value = '  -2.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '552.08'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.97'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.49%  '
$ws.Range("E33").Value = '  -3.45%  '
$ws.Range("E34").Value = '  -1.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.129'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("E37").Value = '  -5.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.40'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.29%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '154.13'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.372'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.27'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.82'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.66%  '
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.53'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -4.70%  '
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.19'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.66%  '
$ws.Range("D47").Value = '0.0₆0296'
$ws.Range("E47").Value = '  -6.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.585'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '152.63'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.63%  '
$ws.Range("E50").Value = '  -2.54%  '
$ws.Range("E51").Value = '  -2.70%  '
